# Update dashboards - 2025-12-04
# Refreshes the "Latest Period" dates and rolling Present/Lag1..Lag4
# observation columns (Q:U) for several rows of the Macro Dashboard,
# as new data points rolled in and older ones shifted down the lag
# window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 (UI Initial Claims / ICSA) ---------------------------------
# Date cell N13 also gets promoted to the "highlighted" date style, used
# elsewhere (e.g. N29) to flag the most recently refreshed series. Copy
# that formatting over before writing the new values.
$ws.Range("N29").Copy()
$ws.Range("N13").PasteSpecial(-4122)
$ws.Range("N13").Value = 45985
$ws.Range("Q13").Value = 191000
$ws.Range("R13").Value = 218000
$ws.Range("S13").Value = 222000
$ws.Range("T13").Value = 228000
$ws.Range("U13").Value = 229000

# --- Row 14 (UI Continuing Claims / CCSA) ------------------------------
$ws.Range("N29").Copy()
$ws.Range("N14").PasteSpecial(-4122)
$ws.Range("N14").Value = 45978
$ws.Range("Q14").Value = 1939000
$ws.Range("R14").Value = 1943000
$ws.Range("S14").Value = 1953000
$ws.Range("T14").Value = 1946000
$ws.Range("U14").Value = 1964000

# --- Row 29 (5yr, 5yr Forward / T5YIFR) --------------------------------
$ws.Range("N29").Value = 45994
$ws.Range("Q29").Value = 2.19
$ws.Range("R29").Value = 2.18
$ws.Range("S29").Value = 2.17
$ws.Range("U29").Value = 0

# --- Row 30 (10yr TIPS / T10YIE) ---------------------------------------
$ws.Range("N30").Value = 45994
$ws.Range("S30").Value = 2.24
$ws.Range("U30").Value = 0

# --- Row 47 (FFR / DFF) -------------------------------------------------
$ws.Range("N47").Value = 45993
$ws.Range("U47").Value = 3.89

# --- Row 48 (2y UST / DGS2) --------------------------------------------
$ws.Range("N48").Value = 45993
$ws.Range("Q48").Value = 3.51
$ws.Range("R48").Value = 3.54
$ws.Range("T48").Value = 0
$ws.Range("U48").Value = 3.47

# --- Row 49 (5y UST / DGS5) --------------------------------------------
$ws.Range("N49").Value = 45993
$ws.Range("Q49").Value = 3.66
$ws.Range("R49").Value = 3.67
$ws.Range("T49").Value = 0
$ws.Range("U49").Value = 3.59

# --- Row 50 (10y UST / DGS10) ------------------------------------------
$ws.Range("N50").Value = 45993
$ws.Range("R50").Value = 4.09
$ws.Range("T50").Value = 0
$ws.Range("U50").Value = 4.02

# --- Row 52 (BAA / DBAA) ------------------------------------------------
$ws.Range("N52").Value = 45993
$ws.Range("Q52").Value = 5.85
$ws.Range("R52").Value = 5.87
$ws.Range("T52").Value = 0
$ws.Range("U52").Value = 5.8
